$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "90.534.15"
$ws.Range("E2").Value2 = "  -0.15%  "
$ws.Range("D3").Value2 = "3.065.71"
$ws.Range("E3").Value2 = "  -1.56%  "
$ws.Range("E4").Value2 = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "243.51"
$ws.Range("E5").Value2 = "  +2.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "617.03"
$ws.Range("E6").Value2 = "  -2.02%  "
$ws.Range("E7").Value2 = "  +7.65%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.364"
$ws.Range("E8").Value2 = "  +1.01%  "
$ws.Range("E9").Value2 = "  +0.00%  "
$ws.Range("D10").Value2 = "3.069.84"
$ws.Range("E10").Value2 = "  -1.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.738"
$ws.Range("E11").Value2 = "  +3.33%  "
$ws.Range("E12").Value2 = "  +2.98%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "0.0000247"
$ws.Range("E13").Value2 = "  +0.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "34.86"
$ws.Range("E14").Value2 = "  -4.32%  "
$ws.Range("D15").Value2 = "90.929.50"
$ws.Range("E15").Value2 = "  +0.38%  "
$ws.Range("E16").Value2 = "  -1.00%  "
$ws.Range("D17").Value2 = "3.652.11"
$ws.Range("E17").Value2 = "  -1.02%  "
$ws.Range("D18").Value2 = "3.118.15"
$ws.Range("E18").Value2 = "  -1.29%  "
$ws.Range("E19").Value2 = "  -1.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "14.40"
$ws.Range("E20").Value2 = "  +2.12%  "
$ws.Range("E21").Value2 = "  +1.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "5.74"
$ws.Range("E22").Value2 = "  +3.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "439.05"
$ws.Range("E23").Value2 = "  -0.54%  "
$ws.Range("E24").Value2 = "  +1.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "90.76"
$ws.Range("E25").Value2 = "  +3.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "5.57"
$ws.Range("E26").Value2 = "  -5.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "11.74"
$ws.Range("E27").Value2 = "  -5.17%  "
$ws.Range("E28").Value2 = "  -0.96%  "
$ws.Range("E29").Value2 = "  +0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "0.250"
$ws.Range("E30").Value2 = "  +28.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "0.181"
$ws.Range("E31").Value2 = "  +13.74%  "
$ws.Range("B32").Value2 = "Binance-PegBSC-USD"
$ws.Range("C32").Value2 = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "1.01"
$ws.Range("E32").Value2 = "  +13.68%  "
$ws.Range("B33").Value2 = "InternetComputer(DFINITY)"
$ws.Range("C33").Value2 = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "9.10"
$ws.Range("E33").Value2 = "  -4.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "0.166"
$ws.Range("E34").Value2 = "  +12.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "0.111"
$ws.Range("E35").Value2 = "  +31.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "7.71"
$ws.Range("E36").Value2 = "  +9.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "26.26"
$ws.Range("E37").Value2 = "  +0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "4.18"
$ws.Range("E38").Value2 = "  +28.91%  "
$ws.Range("E39").Value2 = "  -0.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "490.27"
$ws.Range("E40").Value2 = "  -3.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "3.60"
$ws.Range("E41").Value2 = "  -4.71%  "
$ws.Range("E42").Value2 = "  +0.40%  "
$ws.Range("E43").Value2 = "  +1.23%  "
$ws.Range("E44").Value2 = "  -0.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "153.90"
$ws.Range("E46").Value2 = "  +1.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "1.87"
$ws.Range("E47").Value2 = "  -0.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "0.681"
$ws.Range("E48").Value2 = "  -0.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "4.43"
$ws.Range("E49").Value2 = "  -0.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "1.32"
$ws.Range("E50").Value2 = "  -0.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "44.02"
$ws.Range("E51").Value2 = "  -2.43%  "
